# Add more KPI rows (Profits & Customer Satisfaction, Operational KPI) to the
# "Aerospace" worksheet of the "Domains and KPI's" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aerospace")

# --- Row 21: Profits and Customer Satisfaction (category) -----------------
$ws.Range("A21").Value = "Profits and Customer Satisfaction"
$ws.Range("B21").Value = "Revenue per Kilometer"
$ws.Range("C21").Value = "Volume of passanger carried by aircraft"

# --- Row 22: Passanger Yield ------------------------------------------------
$ws.Range("B22").Value = "Passanger Yield"
$ws.Range("C22").Value = "Passenger yield measures the average revenue received per passenger kilometre flown, calculated by dividing total passenger revenue by the total number of revenue passenger kilometres (RPK)."
$ws.Rows.Item(22).RowHeight = 30

# --- Row 23: Revenue Per Available Seat Kilometer --------------------------
$ws.Range("B23").Value = "Revenue Per Available Seat Kilometer"
$ws.Range("C23").Value = "Revenue per Available Seat Kilometer (RASK) measures the airline’s revenue generated per kilometre for each seat available for sale. It’s calculated by dividing total passenger revenue by available seat kilometres (ASK)."
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24: Revenue per kilometer ------------------------------------------
$ws.Range("B24").Value = "Revenue per kilometer"
$ws.Range("C24").Value = ". It indicates how effectively an airline monetises the distance passengers travel, reflecting pricing strategies and revenue generation efficiency."
$ws.Rows.Item(24).RowHeight = 30

# --- Row 25: Cost Per Available Seat Kilometer ------------------------------
$ws.Range("B25").Value = "Cost Per Available Seat Kilometer"
$ws.Range("C25").Value = "Cost per Available Seat Kilometer (CASK) measures the operating expenses incurred by an airline for each available seat kilometer (ASK), calculated by dividing total operating expenses by ASK."
$ws.Rows.Item(25).RowHeight = 30

# --- Row 26: Break Even Loader ----------------------------------------------
$ws.Range("B26").Value = "Break Even Loader"
$ws.Range("C26").Value = "Break-even Load Factor (BLF) is the passenger load factor at which an airline covers all its costs with revenue, resulting in neither profit nor loss. It helps determine the minimum occupancy level required for flights to break even financially."
$ws.Rows.Item(26).RowHeight = 45

# --- Row 27: Operational KPI (category) / On Time Performance --------------
$ws.Range("A27").Value = "Operational KPI"
$ws.Range("B27").Value = "On Time Performance"
$ws.Range("C27").Value = "On-Time Performance (OTP) measures the percentage of flights that depart and arrive within a specified timeframe of their scheduled times,"
$ws.Rows.Item(27).RowHeight = 30

# --- Row 28: Baggage Handling ------------------------------------------------
$ws.Range("B28").Value = "Baggage Handling"
$ws.Range("C28").Value = "Baggage Handling Performance measures the rate at which baggage is mishandled, encompassing lost, delayed, or damaged luggage incidents. It directly reflects the effectiveness of an airline’s baggage handling processes and operational efficiency."
$ws.Rows.Item(28).RowHeight = 45

# --- Row 29: Customer Satisfaction Index -------------------------------------
$ws.Range("B29").Value = "Customer Satisfaction Index"
$ws.Range("C29").Value = "Customer Satisfaction Index (CSI) is a metric that quantifies passengers’ satisfaction levels with various aspects of airline services, including check-in, onboard experience, and customer service interactions."
$ws.Rows.Item(29).RowHeight = 30

# --- Row 30: Net Promoter Score -----------------------------------------------
$ws.Range("B30").Value = "Net Promoter Score"
$ws.Range("C30").Value = "Net Promoter Score (NPS) assesses customer loyalty and satisfaction based on the likelihood of passengers recommending the airline to others. It categorises customers as promoters (loyal enthusiasts), detractors (unhappy customers), or passives (neutral customers)."
$ws.Rows.Item(30).RowHeight = 45

# Apply the same wrap-text styling used by the rest of column C to the new
# description cells (row 31 is a trailing styled-but-empty cell).
$ws.Range("C21:C31").WrapText = $true

# Update the view: drop the scrolled topLeftCell and select the newly added
# description range, matching the author's final selection state.
$ws.Range("C22:C31").Select()
